$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.485.91"
$ws.Range("E2").Value = "  -1.00%  "
$ws.Range("E3").Value = "  -1.60%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Formula = "'239.15"
$ws.Range("E5").Value = "  -1.25%  "
$ws.Range("D6").Formula = "'0.9998"
$ws.Range("E6").Value = "  -0.07%  "
$ws.Range("D7").Formula = "'0.4778"
$ws.Range("E7").Value = "  -2.25%  "
$ws.Range("D8").Formula = "'0.2834"
$ws.Range("E8").Value = "  -3.47%  "
$ws.Range("D9").Formula = "'0.06704"
$ws.Range("E9").Value = "  -3.07%  "
$ws.Range("D10").Formula = "'18.74"
$ws.Range("E10").Value = "  -3.62%  "
$ws.Range("D11").Formula = "'101.29"
$ws.Range("E11").Value = "  -4.60%  "
$ws.Range("D12").Value = "1.917.91"
$ws.Range("E12").Value = "  -1.49%  "
$ws.Range("D13").Formula = "'0.07673"
$ws.Range("D14").Formula = "'5.203"
$ws.Range("E14").Value = "  -2.73%  "
$ws.Range("D15").Formula = "'0.6687"
$ws.Range("E15").Value = "  -4.53%  "
$ws.Range("D16").Value = "30.502.08"
$ws.Range("E16").Value = "  -0.95%  "
$ws.Range("D17").Formula = "'259.28"
$ws.Range("E17").Value = "  -6.29%  "
$ws.Range("D18").Formula = "'1.000"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("D19").Formula = "'0.000007466"
$ws.Range("E19").Value = "  -3.52%  "
$ws.Range("E20").Value = "  -3.79%  "
$ws.Range("D21").Formula = "'5.388"
$ws.Range("E21").Value = "  -0.98%  "
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").Formula = "'6.280"
$ws.Range("E23").Value = "  -3.81%  "
$ws.Range("D24").Formula = "'9.352"
$ws.Range("E24").Value = "  -3.86%  "
$ws.Range("D25").Formula = "'167.40"
$ws.Range("E25").Value = "  -0.75%  "
$ws.Range("D26").Formula = "'19.14"
$ws.Range("E26").Value = "  -2.49%  "
$ws.Range("D27").Formula = "'2.057"
$ws.Range("E27").Value = "  -4.99%  "
$ws.Range("D28").Formula = "'4.803"
$ws.Range("E28").Value = "  +5.64%  "
$ws.Range("D29").Formula = "'1.386"
$ws.Range("E29").Value = "  -0.58%  "
$ws.Range("E30").Value = "  -3.81%  "
$ws.Range("D31").Formula = "'1.509"
$ws.Range("E31").Value = "  -2.70%  "
$ws.Range("D32").Formula = "'4.256"
$ws.Range("E32").Value = "  -2.54%  "
$ws.Range("D33").Formula = "'0.04713"
$ws.Range("E33").Value = "  -2.96%  "
$ws.Range("D34").Formula = "'0.7250"
$ws.Range("E34").Value = "  -3.39%  "
$ws.Range("E35").Value = "  -4.50%  "
$ws.Range("D36").Formula = "'2.704"
$ws.Range("E36").Value = "  -0.77%  "
$ws.Range("D37").Formula = "'0.01914"
$ws.Range("E37").Value = "  -3.80%  "
$ws.Range("D38").Formula = "'2.617"
$ws.Range("E38").Value = "  -1.98%  "
$ws.Range("D39").Formula = "'6.248"
$ws.Range("E39").Value = "  -3.15%  "
$ws.Range("D40").Formula = "'74.82"
$ws.Range("E40").Value = "  -4.04%  "
$ws.Range("D41").Formula = "'1.967"
$ws.Range("E41").Value = "  -6.72%  "
$ws.Range("D42").Formula = "'0.8612"
$ws.Range("E42").Value = "  -4.85%  "
$ws.Range("D43").Formula = "'105.58"
$ws.Range("E43").Value = "  -2.84%  "
$ws.Range("D44").Formula = "'0.4244"
$ws.Range("E44").Value = "  -3.73%  "
$ws.Range("D45").Formula = "'0.9992"
$ws.Range("E45").Value = "  +0.14%  "
$ws.Range("D46").Formula = "'7.373"
$ws.Range("E46").Value = "  -4.84%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Formula = "'0.1199"
$ws.Range("E47").Value = "  -3.71%  "
$ws.Range("B48").Value = "Maker"
$ws.Range("C48").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D48").Formula = "'919.14"
$ws.Range("E48").Value = "  -7.11%  "
$ws.Range("D49").Formula = "'34.71"
$ws.Range("E49").Value = "  -3.23%  "
$ws.Range("D50").Formula = "'8.754"
$ws.Range("E50").Value = "  -5.50%  "
$ws.Range("D51").Formula = "'0.05741"
$ws.Range("E51").Value = "  +0.18%  "
